# Update "ABP Studio - The Missing Tool for Dotnet Developers.pptx"
# Slide 16: shrink the "abp.io/pricing" textbox and add a red
# "LAST 2 DAYS..!" callout with a Wingdings glyph, wired up with a
# second-click fade-in animation (alongside the existing TextBox 9).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)

$EMU = 12700  # EMU per point
# The COM->OOXML EMU round-trip truncates at float32 precision, which can
# land one EMU short of the target for some point values; a tiny epsilon
# (well below 1/12700 pt) nudges it back onto the exact EMU without
# perceptibly moving/resizing anything.
$EPS = 0.00002

# --- 1. Shrink the existing "TextBox 9" (https://abp.io/pricing) -----------
$pricingBox = $s.Shapes.Item(5)   # id=10, name="TextBox 9"
$pricingBox.Width = (2352047 / $EMU) + $EPS

# --- 2. Add the new "LAST 2 DAYS..!" textbox --------------------------------
$newBox = $s.Shapes.AddTextbox(
    1, `
    (5755505 / $EMU) + $EPS, `
    (5805302 / $EMU) + $EPS, `
    (4803244 / $EMU) + $EPS, `
    (369332 / $EMU) + $EPS)

$newBox.Fill.Visible = $false
$newBox.TextFrame.WordWrap = -1
$newBox.TextFrame.AutoSize = 1

$tr = $newBox.TextFrame.TextRange
$tr.Text = "LAST 2 DAYS..! " + [char]0xF04A
$tr.Font.Bold = $true
$tr.Font.Color.RGB = 255
$tr.ParagraphFormat.Alignment = 3

# Last character is the Wingdings glyph.
$tr.Characters($tr.Length, 1).Font.Name = "Wingdings"

# --- 3. Animation: second click reveals TextBox 9 + the new textbox --------
$tl = $s.TimeLine
$tl.MainSequence.AddEffect($pricingBox, 10, 0, 1) | Out-Null
$tl.MainSequence.AddEffect($newBox, 10, 0, 2) | Out-Null
